$wb = $excel.ActiveWorkbook

# --- Update the sample patient record on the "appointment" sheet -----------
$ws = $wb.Worksheets.Item("appointment")
$ws.Range("A2").Value = "Jan Levinson"
$ws.Range("C2").Value = "9034500777"

# --- De-duplicate the redundant Arial-10-black font/style on
#     "drugs_frequency" -----------------------------------------------------
# F1, G1, E2, F2 and G2 were using a style that pointed at a font entry that
# was identical (in look) to the workbook's normal font. Re-applying the
# same font folds these cells back onto the shared/default style, so the
# duplicate font+style pair is no longer referenced.
$ws2 = $wb.Worksheets.Item("drugs_frequency")
$ws2.Range("F1:G1").Font.Name = "Arial"
$ws2.Range("E2:G2").Font.Name = "Arial"

# --- Switch the active sheet/tab from "drugs_frequency" to "appointment" ---
$ws.Activate()
